# "hours update and TAR update"
# Adds two new status-report rows (1/26/2010) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force the new date cells to be stored as plain text (matching the rest of
# column A, which holds dates as text rather than Excel date serials), then
# strip the temporary number format so no stray style sticks around.
$ws.Range("A26:A27").NumberFormat = "@"
$ws.Range("A26").Value = "1/26/2010"
$ws.Range("A27").Value = "1/26/2010"
$ws.Range("A26:A27").ClearFormats()

$ws.Range("B26").Value = 2
$ws.Range("C26").Value = "Group Meeting"

$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "Weekly Meeting"

# Update the active selection to match the new end of the data.
$ws.Range("C28").Select()
